$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear old content that is no longer used (column A rows that lose their text) ---
$ws.Range("A1").ClearContents()
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()
$ws.Range("A7").ClearContents()
$ws.Range("A8").ClearContents()
$ws.Range("A9").ClearContents()
$ws.Range("A10").ClearContents()

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 28.592447916666668
$ws.Columns("B").ColumnWidth = 61.166666666666664

# --- Header row (bold) ---
$ws.Range("B1").Value2 = "Description"
$ws.Range("C1").Value2 = "Lot"
$ws.Range("D1").Value2 = "Etat"
$ws.Range("E1").Value2 = "Priorité"
$ws.Range("B1:E1").Font.Bold = $true

# --- Section: Mise en place de la structure ---
$ws.Range("A2").Value2 = "Mise en place de la structure :"
$ws.Range("B2:E2").Font.Bold = $true

$ws.Range("B3").Value2 = "Implémenter les modules"
$ws.Range("C3").Value2 = 1

$ws.Range("B4").Value2 = "Implémenter les étudiants"
$ws.Range("C4").Value2 = 1

$ws.Range("B5").Value2 = "Implémenter les groupes de modules"
$ws.Range("C5").Value2 = 1

# --- Section: Gestion du csv ---
$ws.Range("A6").Value2 = "Gestion du csv :"
$ws.Range("B6:E6").Font.Bold = $true

$ws.Range("B7").Value2 = "importer les données à partir d'un fichier csv"
$ws.Range("B8").Value2 = "Convertir des fichiers Excel en csv"
$ws.Range("B9").Value2 = "Traiter les données du csv"
$ws.Range("C9").Value2 = 1

# --- Section: Interface Homme-Machine ---
$ws.Range("A11").Value2 = "Interface Homme-Machine :"
$ws.Range("B12").Value2 = "Liste d'étudiants à choisir pour mettre un avis"
$ws.Range("B13").Value2 = "Générer des graphiques "

# --- Section: Gestion de module ---
$ws.Range("A15").Value2 = "Gestion de module :"
$ws.Range("B16").Value2 = "Créer des groupes de modules"
$ws.Range("B17").Value2 = "Gérer des groupes de modules"
$ws.Range("B18").Value2 = "Supprimer des groupes de modules"

# --- Section: Gestion des avis ---
$ws.Range("A20").Value2 = "Gestion des avis : "
$ws.Range("B21").Value2 = "Mettre des avis sur les etudiants"
$ws.Range("B22").Value2 = "Permettre de générer un avis général"

# --- Section: Gestion des données ---
$ws.Range("A24").Value2 = "Gestion des données :"
$ws.Range("B25").Value2 = "Ecriture dans un fichier des données pour enregistrer"
$ws.Range("B26").Value2 = "Génération de pdf"

# --- Underline every cell in column A used by the backlog (section titles column) ---
$ws.Range("A1:A24").Font.Underline = $true

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
